$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.291.57"
$ws.Range("E2").Value = "  -4.00%  "
$ws.Range("D3").Value = "3.026.70"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'538.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.71%  "
$ws.Range("D6").Value = "'132.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -10.62%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.019.31"
$ws.Range("E8").Value = "  -3.53%  "
$ws.Range("D9").Value = "'0.482"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("D10").Value = "'6.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.05%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("D12").Value = "'0.453"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").Value = "'34.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.00%  "
$ws.Range("D14").Value = "'0.0000211"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.89%  "
$ws.Range("D15").Value = "3.509.32"
$ws.Range("E15").Value = "  -3.65%  "
$ws.Range("D16").Value = "62.265.13"
$ws.Range("E16").Value = "  -4.12%  "
$ws.Range("E17").Value = "  -2.53%  "
$ws.Range("D18").Value = "3.026.96"
$ws.Range("E18").Value = "  -3.54%  "
$ws.Range("D19").Value = "'6.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("D20").Value = "'474.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.87%  "
$ws.Range("D21").Value = "'13.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("D22").Value = "'0.686"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").Value = "'6.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.13%  "
$ws.Range("D24").Value = "'76.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.73%  "
$ws.Range("D25").Value = "'11.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.81%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'8.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.92%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.11%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.49%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'25.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("D33").Value = "'59.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.12%  "
$ws.Range("D34").Value = "'2.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.75%  "
$ws.Range("D35").Value = "'507.80"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.65%  "
$ws.Range("D36").Value = "'5.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.31%  "
$ws.Range("D37").Value = "'5.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.13%  "
$ws.Range("D38").Value = "'0.0393"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.70%  "
$ws.Range("D39").Value = "3.032.47"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "'0.0776"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.81%  "
$ws.Range("D41").Value = "'0.116"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.04%  "
$ws.Range("D42").Value = "'7.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("D43").Value = "'2.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.37%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'0.248"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.93%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'1.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.90%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'119.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").Value = "'23.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.68%  "
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("D50").Value = "'2.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +63.37%  "
$ws.Range("D51").Value = "0.0₃0483"
$ws.Range("E51").Value = "  -7.66%  "
